{"js": "// administrative_review.docx edit:\n//  1. Drop the stray \"_GoBack\" bookmark that sat in the\n//     `dor_cse_notice_arrears_wrong` checkbox line.\n//  2. Extend the `account_same` condition to\n//     `account_same is defined  and benefits_used` (shared-account-only\n//     fix referenced in the commit message).\n//  3. Re-insert a \"_GoBack\" bookmark inside the `dor_take_ward_money`\n//     checkbox line (landing right before the final \"y\" of the second\n//     `dor_take_ward_money` occurrence there), mirroring where Word\n//     parked the cursor after the last edit.\n\nconst doc = context.document;\nconst body = doc.body;\n\n// 1) Remove the old \"_GoBack\" bookmark (it was next to\n//    \"{% if dor_cse_notice_arrears_wrong %} [X]{% else %}\").\nconst oldBookmark = doc.getBookmarkRangeOrNullObject(\"_GoBack\");\noldBookmark.load(\"isNullObject\");\nawait context.sync();\nif (!oldBookmark.isNullObject) {\n  doc.deleteBookmark(\"_GoBack\");\n  await context.sync();\n}\n\n// 2) \"{% if account_same %}\" -> \"{% if account_same is  defined  and benefits_used  %}\"\nconst accountSameMatches = body.search(\"account_same\", { matchCase: true });\naccountSameMatches.load(\"items\");\nawait context.sync();\nif (accountSameMatches.items.length > 0) {\n  const accountSameRange = accountSameMatches.items[0];\n  accountSameRange.insertText(\n    \" is defined  and benefits_used\",\n    Word.InsertLocation.end\n  );\n  await context.sync();\n}\n\n// 3) Re-insert \"_GoBack\" just before the trailing \"y\" of the second\n//    \"dor_take_ward_money\" in \"{% if  dor_take_ward_money is defined and\n//    dor_take_ward_money %}[X]{% else %}[  ]{% endif %}\".\nconst wardMoneyMatches = body.search(\"dor_take_ward_mone\", { matchCase: true });\nwardMoneyMatches.load(\"items\");\nawait context.sync();\nif (wardMoneyMatches.items.length >= 2) {\n  const secondOccurrence = wardMoneyMatches.items[1];\n  const insertionPoint = secondOccurrence.getRange(Word.RangeLocation.end);\n  insertionPoint.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# administrative_review.docx edit:\n#  1. Drop the stray \"_GoBack\" bookmark that sat in the\n#     `dor_cse_notice_arrears_wrong` checkbox line.\n#  2. Extend the `account_same` condition to\n#     `account_same is defined  and benefits_used` (shared-account-only\n#     fix referenced in the commit message).\n#  3. Re-insert a \"_GoBack\" bookmark inside the `dor_take_ward_money`\n#     checkbox line (landing right before the final \"y\" of the second\n#     `dor_take_ward_money` occurrence there), mirroring where Word\n#     parked the cursor after the last edit.\n\n$d = $word.ActiveDocument\n\n# 1) Remove the old \"_GoBack\" bookmark.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# 2) \"{% if account_same %}\" -> \"{% if account_same is  defined  and benefits_used  %}\"\n$findRange = $d.Content\n$find = $findRange.Find\n$find.Text = \"account_same\"\n$find.MatchCase = $true\nif ($find.Execute()) {\n    $findRange.Collapse(0)  # wdCollapseEnd\n    $findRange.InsertAfter(\" is defined  and benefits_used\")\n}\n\n# 3) Re-insert \"_GoBack\" just before the trailing \"y\" of the second\n#    \"dor_take_ward_money\" in \"{% if  dor_take_ward_money is defined and\n#    dor_take_ward_money %}[X]{% else %}[  ]{% endif %}\".\n$wardRange = $d.Content\n$wardFind = $wardRange.Find\n$wardFind.Text = \"dor_take_ward_mone\"\n$wardFind.MatchCase = $true\nif ($wardFind.Execute()) {\n    if ($wardFind.Execute()) {\n        $wardRange.Collapse(0)  # wdCollapseEnd\n        $d.Bookmarks.Add(\"_GoBack\", $wardRange) | Out-Null\n    }\n}\n"}
